# Apply coin price / label updates for the 2022-12-20 GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'247.41"
$ws.Range("D3").Value = "'22.04"
$ws.Range("D4").Value = "'5.340"
$ws.Range("D5").Value = "'0.05637"
$ws.Range("D6").Value = "'3.426"
$ws.Range("D7").Value = "'6.363"
$ws.Range("D8").Value = "'0.8186"
$ws.Range("D9").Value = "'0.9372"
$ws.Range("D10").Value = "'0.1440"
$ws.Range("D11").Value = "'0.07509"
$ws.Range("D12").Value = "'0.03245"
$ws.Range("D13").Value = "'0.03083"
$ws.Range("D14").Value = "'0.09304"
$ws.Range("D15").Value = "'3.567"
$ws.Range("D16").Value = "'0.001600"
$ws.Range("D17").Value = "'0.04737"
$ws.Range("D18").Value = "'0.0005777"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("D19").Value = "'0.006319"
$ws.Range("D20").Value = "'0.005064"
$ws.Range("D21").Value = "'0.001035"
$ws.Range("D22").Value = "'0.0001500"
$ws.Range("D23").Value = "'3.754"
$ws.Range("D24").Value = "'2.150"
$ws.Range("D25").Value = "'0.3306"
$ws.Range("D26").Value = "'0.1318"
$ws.Range("D28").Value = "'0.0002998"
$ws.Range("D40").Value = "'0.03955"
$ws.Range("D41").Value = "'0.006990"
$ws.Range("D42").Value = "'0.1064"
$ws.Range("D43").Value = "'0.003400"
$ws.Range("D44").Value = "'0.008661"
$ws.Range("D45").Value = "'0.00005576"
$ws.Range("D47").Value = "'0.0005497"
$ws.Range("E47").Value = "46ACDXExchangeACXT"
$ws.Range("D48").Value = "'0.7794"
$ws.Range("D49").Value = "'0.1778"
$ws.Range("E49").Value = "48BOLOBOLOBestin24h"
$ws.Range("D50").Value = "'0.00002100"
